$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows 20-21 need the same date number format as the rest of column D
$ws.Cells.Item(20, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(21, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 11
$ws.Cells.Item(11, 1).Value = 2
$ws.Cells.Item(11, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(11, 3).Value = "Coquimbo"
$ws.Cells.Item(11, 4).Value = 44566
$ws.Cells.Item(11, 5).Value = 4
$ws.Cells.Item(11, 6).Value = "Fruta"
$ws.Cells.Item(11, 7).Value = 100103
$ws.Cells.Item(11, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(11, 9).Value = 100103003
$ws.Cells.Item(11, 10).Value = "Damasco"
$ws.Cells.Item(11, 11).Value = "Modesto"
$ws.Cells.Item(11, 12).Value = "Especial"
$ws.Cells.Item(11, 13).Value = 100
$ws.Cells.Item(11, 14).Value = 23000
$ws.Cells.Item(11, 15).Value = 24000
$ws.Cells.Item(11, 16).Value = 23500
$ws.Cells.Item(11, 17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(11, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(11, 19).Value = 1306
$ws.Cells.Item(11, 20).Value = 18

# Row 12
$ws.Cells.Item(12, 1).Value = 2
$ws.Cells.Item(12, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(12, 3).Value = "Coquimbo"
$ws.Cells.Item(12, 4).Value = 44566
$ws.Cells.Item(12, 5).Value = 4
$ws.Cells.Item(12, 6).Value = "Fruta"
$ws.Cells.Item(12, 7).Value = 100103
$ws.Cells.Item(12, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(12, 9).Value = 100103003
$ws.Cells.Item(12, 10).Value = "Damasco"
$ws.Cells.Item(12, 11).Value = "Modesto"
$ws.Cells.Item(12, 12).Value = "Primera"
$ws.Cells.Item(12, 13).Value = 160
$ws.Cells.Item(12, 14).Value = 21000
$ws.Cells.Item(12, 15).Value = 22000
$ws.Cells.Item(12, 16).Value = 21500
$ws.Cells.Item(12, 17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(12, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(12, 19).Value = 1194
$ws.Cells.Item(12, 20).Value = 18

# Row 13
$ws.Cells.Item(13, 1).Value = 2
$ws.Cells.Item(13, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(13, 3).Value = "Coquimbo"
$ws.Cells.Item(13, 4).Value = 44553
$ws.Cells.Item(13, 5).Value = 4
$ws.Cells.Item(13, 6).Value = "Fruta"
$ws.Cells.Item(13, 7).Value = 100103
$ws.Cells.Item(13, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(13, 9).Value = 100103003
$ws.Cells.Item(13, 10).Value = "Damasco"
$ws.Cells.Item(13, 11).Value = "Modesto"
$ws.Cells.Item(13, 12).Value = "Especial"
$ws.Cells.Item(13, 13).Value = 360
$ws.Cells.Item(13, 14).Value = 23000
$ws.Cells.Item(13, 15).Value = 24000
$ws.Cells.Item(13, 16).Value = 23500
$ws.Cells.Item(13, 17).Value = "`$/caja 16 kilos"
$ws.Cells.Item(13, 18).Value = "Región Metropolitana"
$ws.Cells.Item(13, 19).Value = 1469
$ws.Cells.Item(13, 20).Value = 16

# Row 14
$ws.Cells.Item(14, 1).Value = 2
$ws.Cells.Item(14, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(14, 3).Value = "Coquimbo"
$ws.Cells.Item(14, 4).Value = 44553
$ws.Cells.Item(14, 5).Value = 4
$ws.Cells.Item(14, 6).Value = "Fruta"
$ws.Cells.Item(14, 7).Value = 100103
$ws.Cells.Item(14, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(14, 9).Value = 100103003
$ws.Cells.Item(14, 10).Value = "Damasco"
$ws.Cells.Item(14, 11).Value = "Modesto"
$ws.Cells.Item(14, 12).Value = "Primera"
$ws.Cells.Item(14, 13).Value = 300
$ws.Cells.Item(14, 14).Value = 21000
$ws.Cells.Item(14, 15).Value = 22000
$ws.Cells.Item(14, 16).Value = 21500
$ws.Cells.Item(14, 17).Value = "`$/caja 16 kilos"
$ws.Cells.Item(14, 18).Value = "Región Metropolitana"
$ws.Cells.Item(14, 19).Value = 1344
$ws.Cells.Item(14, 20).Value = 16

# Row 15
$ws.Cells.Item(15, 1).Value = 2
$ws.Cells.Item(15, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(15, 3).Value = "Coquimbo"
$ws.Cells.Item(15, 4).Value = 44553
$ws.Cells.Item(15, 5).Value = 4
$ws.Cells.Item(15, 6).Value = "Fruta"
$ws.Cells.Item(15, 7).Value = 100103
$ws.Cells.Item(15, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(15, 9).Value = 100103003
$ws.Cells.Item(15, 10).Value = "Damasco"
$ws.Cells.Item(15, 11).Value = "Modesto"
$ws.Cells.Item(15, 12).Value = "Segunda"
$ws.Cells.Item(15, 13).Value = 240
$ws.Cells.Item(15, 14).Value = 17000
$ws.Cells.Item(15, 15).Value = 18000
$ws.Cells.Item(15, 16).Value = 17500
$ws.Cells.Item(15, 17).Value = "`$/caja 16 kilos"
$ws.Cells.Item(15, 18).Value = "Región Metropolitana"
$ws.Cells.Item(15, 19).Value = 1094
$ws.Cells.Item(15, 20).Value = 16

# Row 16
$ws.Cells.Item(16, 1).Value = 2
$ws.Cells.Item(16, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(16, 3).Value = "Coquimbo"
$ws.Cells.Item(16, 4).Value = 44161
$ws.Cells.Item(16, 5).Value = 4
$ws.Cells.Item(16, 6).Value = "Fruta"
$ws.Cells.Item(16, 7).Value = 100103
$ws.Cells.Item(16, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(16, 9).Value = 100103003
$ws.Cells.Item(16, 10).Value = "Damasco"
$ws.Cells.Item(16, 11).Value = "Dina"
$ws.Cells.Item(16, 12).Value = "Primera"
$ws.Cells.Item(16, 13).Value = 300
$ws.Cells.Item(16, 14).Value = 20000
$ws.Cells.Item(16, 15).Value = 20500
$ws.Cells.Item(16, 16).Value = 20250
$ws.Cells.Item(16, 17).Value = "`$/caja 15 kilos"
$ws.Cells.Item(16, 18).Value = "Región Metropolitana"
$ws.Cells.Item(16, 19).Value = 1350
$ws.Cells.Item(16, 20).Value = 15

# Row 17
$ws.Cells.Item(17, 1).Value = 2
$ws.Cells.Item(17, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(17, 3).Value = "Coquimbo"
$ws.Cells.Item(17, 4).Value = 44161
$ws.Cells.Item(17, 5).Value = 4
$ws.Cells.Item(17, 6).Value = "Fruta"
$ws.Cells.Item(17, 7).Value = 100103
$ws.Cells.Item(17, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(17, 9).Value = 100103003
$ws.Cells.Item(17, 10).Value = "Damasco"
$ws.Cells.Item(17, 11).Value = "Dina"
$ws.Cells.Item(17, 12).Value = "Segunda"
$ws.Cells.Item(17, 13).Value = 100
$ws.Cells.Item(17, 14).Value = 18000
$ws.Cells.Item(17, 15).Value = 18500
$ws.Cells.Item(17, 16).Value = 18250
$ws.Cells.Item(17, 17).Value = "`$/caja 15 kilos"
$ws.Cells.Item(17, 18).Value = "Región Metropolitana"
$ws.Cells.Item(17, 19).Value = 1217
$ws.Cells.Item(17, 20).Value = 15

# Row 18
$ws.Cells.Item(18, 1).Value = 2
$ws.Cells.Item(18, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(18, 3).Value = "Coquimbo"
$ws.Cells.Item(18, 4).Value = 44160
$ws.Cells.Item(18, 5).Value = 4
$ws.Cells.Item(18, 6).Value = "Fruta"
$ws.Cells.Item(18, 7).Value = 100103
$ws.Cells.Item(18, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(18, 9).Value = 100103003
$ws.Cells.Item(18, 10).Value = "Damasco"
$ws.Cells.Item(18, 11).Value = "Castle Brite"
$ws.Cells.Item(18, 12).Value = "Primera"
$ws.Cells.Item(18, 13).Value = 240
$ws.Cells.Item(18, 14).Value = 20500
$ws.Cells.Item(18, 15).Value = 21000
$ws.Cells.Item(18, 16).Value = 20750
$ws.Cells.Item(18, 17).Value = "`$/caja 15 kilos"
$ws.Cells.Item(18, 18).Value = "Región Metropolitana"
$ws.Cells.Item(18, 19).Value = 1383
$ws.Cells.Item(18, 20).Value = 15

# Row 19
$ws.Cells.Item(19, 1).Value = 2
$ws.Cells.Item(19, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(19, 3).Value = "Coquimbo"
$ws.Cells.Item(19, 4).Value = 44559
$ws.Cells.Item(19, 5).Value = 4
$ws.Cells.Item(19, 6).Value = "Fruta"
$ws.Cells.Item(19, 7).Value = 100103
$ws.Cells.Item(19, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(19, 9).Value = 100103003
$ws.Cells.Item(19, 10).Value = "Damasco"
$ws.Cells.Item(19, 11).Value = "Modesto"
$ws.Cells.Item(19, 12).Value = "Especial"
$ws.Cells.Item(19, 13).Value = 400
$ws.Cells.Item(19, 14).Value = 25000
$ws.Cells.Item(19, 15).Value = 26000
$ws.Cells.Item(19, 16).Value = 25500
$ws.Cells.Item(19, 17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(19, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(19, 19).Value = 1417
$ws.Cells.Item(19, 20).Value = 18

# Row 20
$ws.Cells.Item(20, 1).Value = 2
$ws.Cells.Item(20, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(20, 3).Value = "Coquimbo"
$ws.Cells.Item(20, 4).Value = 44559
$ws.Cells.Item(20, 5).Value = 4
$ws.Cells.Item(20, 6).Value = "Fruta"
$ws.Cells.Item(20, 7).Value = 100103
$ws.Cells.Item(20, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(20, 9).Value = 100103003
$ws.Cells.Item(20, 10).Value = "Damasco"
$ws.Cells.Item(20, 11).Value = "Modesto"
$ws.Cells.Item(20, 12).Value = "Primera"
$ws.Cells.Item(20, 13).Value = 320
$ws.Cells.Item(20, 14).Value = 22000
$ws.Cells.Item(20, 15).Value = 23000
$ws.Cells.Item(20, 16).Value = 22500
$ws.Cells.Item(20, 17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(20, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(20, 19).Value = 1250
$ws.Cells.Item(20, 20).Value = 18

# Row 21
$ws.Cells.Item(21, 1).Value = 2
$ws.Cells.Item(21, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(21, 3).Value = "Coquimbo"
$ws.Cells.Item(21, 4).Value = 44175
$ws.Cells.Item(21, 5).Value = 4
$ws.Cells.Item(21, 6).Value = "Fruta"
$ws.Cells.Item(21, 7).Value = 100103
$ws.Cells.Item(21, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(21, 9).Value = 100103003
$ws.Cells.Item(21, 10).Value = "Damasco"
$ws.Cells.Item(21, 11).Value = "Castle Brite"
$ws.Cells.Item(21, 12).Value = "Primera"
$ws.Cells.Item(21, 13).Value = 300
$ws.Cells.Item(21, 14).Value = 21000
$ws.Cells.Item(21, 15).Value = 22000
$ws.Cells.Item(21, 16).Value = 21500
$ws.Cells.Item(21, 17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(21, 18).Value = "Región Metropolitana"
$ws.Cells.Item(21, 19).Value = 1194
$ws.Cells.Item(21, 20).Value = 18
